# "fix bugtracker to git ignore"
# Adds a new bug-tracker row (row 9) for:
#   TYPE=feature request, STATUS=open, TARGET=(blank), DATE==TODAY(),
#   component=course/actionCard, SUMMARY=show action card only when its class time

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New row 9 values ---
$ws.Range("A9").Value = "feature request"
$ws.Range("B9").Value = "open"
$ws.Range("E9").Value = "course/actionCard"
$ws.Range("F9").Value = "show action card only when its class time"

# Match the green "open" status fill already used on B4:B8 (reuses the
# existing style instead of fabricating a new one).
$ws.Range("B9").Interior.Color = $ws.Range("B4").Interior.Color

# Long-date format for the new DATE cell, then the live TODAY() formula.
$ws.Range("D9").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"
$ws.Range("D9").Formula = "=TODAY()"

# Sheet view now focuses the new SUMMARY cell instead of the old D9 spot.
$ws.Range("F9").Select()
